$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1615110.8
$ws.Range("J17").Value = 1668925.1
$ws.Range("L17").Value = 5006775.300000001
$ws.Range("N17").Value = -5007111.300000001
$ws.Range("H38").Value = 565.7692
$ws.Range("I38").Value = 59.166668
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 177.500004
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = 194.499996
$ws.Range("N38").Value = -3744
$ws.Range("H74").Value = 5428.5713
$ws.Range("J74").Value = 5333.3335
$ws.Range("L74").Value = 5333.3335
$ws.Range("N74").Value = -7205.3335
$ws.Range("H76").Value = 2184
$ws.Range("I76").Value = 2184
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 2184
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -1869
$ws.Range("N76").ClearContents()
$ws.Range("H77").Value = 5428.5713
$ws.Range("J77").Value = 5333.3335
$ws.Range("L77").Value = 26666.6675
$ws.Range("N77").Value = -36026.6675
$ws.Range("H79").Value = 2184
$ws.Range("I79").Value = 2184
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 2184
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -1092
$ws.Range("N79").ClearContents()
$ws.Range("H137").Value = 2167.1875
$ws.Range("I137").Value = 2052.923
$ws.Range("K137").Value = 6158.768999999999
$ws.Range("M137").Value = -3608.768999999999
$ws.Range("H138").Value = 2107.9495
$ws.Range("J138").Value = 2312.4285
$ws.Range("L138").Value = 6937.2855
$ws.Range("N138").Value = -17217.2855

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 31249996
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 31249996
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H74").Value = 71429640
$ws.Range("I74").Value = 90909670
$ws.Range("K74").Value = 90909670
$ws.Range("M74").Value = -90908796
$ws.Range("H77").Value = 71429640
$ws.Range("I77").Value = 90909670
$ws.Range("K77").Value = 454548350
$ws.Range("M77").Value = -454543982
$ws.Range("H122").Value = 3533.1667
$ws.Range("I122").Value = 2200
$ws.Range("K122").Value = 6600
$ws.Range("M122").Value = -4150

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4164.125
$ws.Range("I20").Value = 9505.333000000001
$ws.Range("K20").Value = 9505.333000000001
$ws.Range("M20").Value = -9258.333000000001
$ws.Range("H134").Value = 6384.4116
$ws.Range("I134").Value = 7042.3335
$ws.Range("J134").Value = 1450
$ws.Range("K134").Value = 21127.0005
$ws.Range("L134").Value = 4350
$ws.Range("M134").Value = -18592.0005
$ws.Range("N134").Value = -9420

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12280.333
$ws.Range("I31").Value = 16696.154
$ws.Range("J31").Value = 5104.625
$ws.Range("K31").Value = 16696.154
$ws.Range("L31").Value = 5104.625
$ws.Range("M31").Value = -16401.154
$ws.Range("N31").Value = -5694.625
$ws.Range("H34").Value = 12280.333
$ws.Range("I34").Value = 16696.154
$ws.Range("J34").Value = 5104.625
$ws.Range("K34").Value = 16696.154
$ws.Range("L34").Value = 5104.625
$ws.Range("M34").Value = -16494.154
$ws.Range("N34").Value = -5508.625
$ws.Range("H62").Value = 250004260
$ws.Range("I62").Value = 250004260
$ws.Range("K62").Value = 250004260
$ws.Range("M62").Value = -250003636
$ws.Range("H65").Value = 250004260
$ws.Range("I65").Value = 250004260
$ws.Range("K65").Value = 1250021300
$ws.Range("M65").Value = -1250018180
$ws.Range("H134").Value = 1191.0513
$ws.Range("I134").Value = 891.78125
$ws.Range("K134").Value = 2675.34375
$ws.Range("M134").Value = -140.34375

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 3889.2856
$ws.Range("J54").Value = 3889.2856
$ws.Range("L54").Value = 11667.8568
$ws.Range("N54").Value = -12785.8568
$ws.Range("H92").Value = 592.46155
$ws.Range("I92").Value = 628.8570999999999
$ws.Range("J92").Value = 550
$ws.Range("K92").Value = 1886.5713
$ws.Range("L92").Value = 1650
$ws.Range("M92").Value = -638.5712999999998
$ws.Range("N92").Value = -4146
$ws.Range("H113").Value = 516.8889
$ws.Range("J113").Value = 574.8125
$ws.Range("L113").Value = 1724.4375
$ws.Range("N113").Value = -6064.4375
$ws.Range("H121").Value = 4183
$ws.Range("I121").Value = 407.5
$ws.Range("J121").Value = 5212.6816
$ws.Range("K121").Value = 1222.5
$ws.Range("L121").Value = 15638.0448
$ws.Range("M121").Value = 87.5
$ws.Range("N121").Value = -18258.0448
$ws.Range("H131").Value = 785.55
$ws.Range("J131").Value = 785.55
$ws.Range("L131").Value = 2356.65
$ws.Range("N131").Value = -12436.65
$ws.Range("H134").Value = 2017.8422
$ws.Range("I134").Value = 1222.3334
$ws.Range("J134").Value = 5001
$ws.Range("K134").Value = 3667.0002
$ws.Range("L134").Value = 15003
$ws.Range("M134").Value = 1402.9998
$ws.Range("N134").Value = -25143
$ws.Range("H137").Value = 25644242
$ws.Range("I137").Value = 904
$ws.Range("J137").Value = 41671330
$ws.Range("K137").Value = 2712
$ws.Range("L137").Value = 125013990
$ws.Range("M137").Value = 2388
$ws.Range("N137").Value = -125024190
$ws.Range("H138").Value = 158991.05
$ws.Range("I138").Value = 1301.875
$ws.Range("J138").Value = 1000000
$ws.Range("K138").Value = 3905.625
$ws.Range("L138").Value = 3000000
$ws.Range("M138").Value = 1234.375
$ws.Range("N138").Value = -3010280

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 10000
$ws.Range("J15").Value = 10000
$ws.Range("L15").Value = 10000
$ws.Range("N15").Value = -10576
$ws.Range("H80").Value = 3308.682
$ws.Range("I80").Value = 2940.5
$ws.Range("J80").Value = 3615.5
$ws.Range("K80").Value = 2940.5
$ws.Range("L80").Value = 3615.5
$ws.Range("M80").Value = -1942.5
$ws.Range("N80").Value = -5611.5
$ws.Range("H81").Value = 10000
$ws.Range("J81").Value = 10000
$ws.Range("L81").Value = 10000
$ws.Range("N81").Value = -11996
$ws.Range("H83").Value = 3308.682
$ws.Range("I83").Value = 2940.5
$ws.Range("J83").Value = 3615.5
$ws.Range("K83").Value = 14702.5
$ws.Range("L83").Value = 18077.5
$ws.Range("M83").Value = -9710.5
$ws.Range("N83").Value = -28061.5
$ws.Range("H84").Value = 10000
$ws.Range("J84").Value = 10000
$ws.Range("L84").Value = 30000
$ws.Range("N84").Value = -39984
$ws.Range("H122").Value = 53334480
$ws.Range("I122").Value = 18519742
$ws.Range("J122").Value = 142858100
$ws.Range("K122").Value = 55559226
$ws.Range("L122").Value = 428574300
$ws.Range("M122").Value = -55556776
$ws.Range("N122").Value = -428579200

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5106.3477
$ws.Range("I40").Value = 3075
$ws.Range("J40").Value = 5823.294
$ws.Range("K40").Value = 3075
$ws.Range("L40").Value = 5823.294
$ws.Range("M40").Value = -2939
$ws.Range("N40").Value = -6095.294
$ws.Range("H122").Value = 1034580.2
$ws.Range("J122").Value = 2811.4285
$ws.Range("L122").Value = 8434.2855
$ws.Range("N122").Value = -13334.2855

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4765
$ws.Range("I62").Value = 4441.6665
$ws.Range("J62").Value = 5250
$ws.Range("K62").Value = 4441.6665
$ws.Range("L62").Value = 5250
$ws.Range("M62").Value = -3817.6665
$ws.Range("N62").Value = -6498
$ws.Range("H65").Value = 4765
$ws.Range("I65").Value = 4441.6665
$ws.Range("J65").Value = 5250
$ws.Range("K65").Value = 22208.3325
$ws.Range("L65").Value = 26250
$ws.Range("M65").Value = -19088.3325
$ws.Range("N65").Value = -32490
$ws.Range("H81").Value = 90910610
$ws.Range("I81").Value = 1670.9
$ws.Range("J81").Value = 1000000000
$ws.Range("K81").Value = 3341.8
$ws.Range("L81").Value = 2000000000
$ws.Range("M81").Value = -2280.8
$ws.Range("N81").Value = -2000002122
$ws.Range("H84").Value = 90910610
$ws.Range("I84").Value = 1670.9
$ws.Range("J84").Value = 1000000000
$ws.Range("K84").Value = 16709
$ws.Range("L84").Value = 10000000000
$ws.Range("M84").Value = -11405
$ws.Range("N84").Value = -10000010608
$ws.Range("H122").Value = 1913.3636
$ws.Range("I122").Value = 1965.625
$ws.Range("J122").Value = 1774
$ws.Range("K122").Value = 5896.875
$ws.Range("L122").Value = 5322
$ws.Range("M122").Value = -3446.875
$ws.Range("N122").Value = -10222
